$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 5).Value = 3
$ws.Cells.Item(2, 7).Value = 61.91334466666666
$ws.Cells.Item(2, 8).Value = 185.740034
$ws.Cells.Item(2, 9).Value = 0.5020829437194911
$ws.Cells.Item(2, 10).Value = 0.5020829437194911
$ws.Cells.Item(2, 11).Value = 3
$ws.Cells.Item(2, 13).Value = 4.358929333333333
$ws.Cells.Item(2, 14).Value = 13.076788
$ws.Cells.Item(2, 15).Value = 0.315486934607403
$ws.Cells.Item(2, 16).Value = 0.315486934607403
$ws.Cells.Item(2, 17).Value = 269.8758941923102
$ws.Cells.Item(2, 18).Value = 2428.883047730792
$ws.Cells.Item(2, 19).Value = 0.1584006088327235
$ws.Cells.Item(2, 20).Value = 0.1584006088327235

$ws.Cells.Item(3, 5).Value = 3
$ws.Cells.Item(3, 7).Value = 61.91334466666666
$ws.Cells.Item(3, 8).Value = 185.740034
$ws.Cells.Item(3, 9).Value = 0.5020829437194911
$ws.Cells.Item(3, 10).Value = 0.5020829437194911
$ws.Cells.Item(3, 11).Value = 3
$ws.Cells.Item(3, 13).Value = 1.781382333333333
$ws.Cells.Item(3, 14).Value = 5.344147
$ws.Cells.Item(3, 15).Value = 0.1289313977653647
$ws.Cells.Item(3, 16).Value = 0.1289313977653647
$ws.Cells.Item(3, 17).Value = 110.2913383867775
$ws.Cells.Item(3, 18).Value = 992.6220454809978
$ws.Cells.Item(3, 19).Value = 0.06473425572790295
$ws.Cells.Item(3, 20).Value = 0.06473425572790295

$ws.Cells.Item(4, 5).Value = 3
$ws.Cells.Item(4, 7).Value = 61.91334466666666
$ws.Cells.Item(4, 8).Value = 185.740034
$ws.Cells.Item(4, 9).Value = 0.5020829437194911
$ws.Cells.Item(4, 10).Value = 0.5020829437194911
$ws.Cells.Item(4, 11).Value = 3
$ws.Cells.Item(4, 13).Value = 7.676201333333334
$ws.Cells.Item(4, 14).Value = 23.028604
$ws.Cells.Item(4, 15).Value = 0.5555816676272323
$ws.Cells.Item(4, 16).Value = 0.5555816676272323
$ws.Cells.Item(4, 17).Value = 475.2592988813929
$ws.Cells.Item(4, 18).Value = 4277.333689932536
$ws.Cells.Item(4, 19).Value = 0.2789480791588647
$ws.Cells.Item(4, 20).Value = 0.2789480791588647

$ws.Cells.Item(5, 5).Value = 3
$ws.Cells.Item(5, 7).Value = 50.09443433333333
$ws.Cells.Item(5, 8).Value = 150.283303
$ws.Cells.Item(5, 9).Value = 0.4062381250674705
$ws.Cells.Item(5, 10).Value = 0.4062381250674706
$ws.Cells.Item(5, 11).Value = 3
$ws.Cells.Item(5, 13).Value = 4.358929333333333
$ws.Cells.Item(5, 14).Value = 13.076788
$ws.Cells.Item(5, 15).Value = 0.315486934607403
$ws.Cells.Item(5, 16).Value = 0.315486934607403
$ws.Cells.Item(5, 17).Value = 218.3580992523071
$ws.Cells.Item(5, 18).Value = 1965.222893270764
$ws.Cells.Item(5, 19).Value = 0.128162820798195
$ws.Cells.Item(5, 20).Value = 0.1281628207981951

$ws.Cells.Item(6, 5).Value = 3
$ws.Cells.Item(6, 7).Value = 50.09443433333333
$ws.Cells.Item(6, 8).Value = 150.283303
$ws.Cells.Item(6, 9).Value = 0.4062381250674705
$ws.Cells.Item(6, 10).Value = 0.4062381250674706
$ws.Cells.Item(6, 11).Value = 3
$ws.Cells.Item(6, 13).Value = 1.781382333333333
$ws.Cells.Item(6, 14).Value = 5.344147
$ws.Cells.Item(6, 15).Value = 0.1289313977653647
$ws.Cells.Item(6, 16).Value = 0.1289313977653647
$ws.Cells.Item(6, 17).Value = 89.23734031972677
$ws.Cells.Item(6, 18).Value = 803.1360628775409
$ws.Cells.Item(6, 19).Value = 0.05237684929053003
$ws.Cells.Item(6, 20).Value = 0.05237684929053004

$ws.Cells.Item(7, 5).Value = 3
$ws.Cells.Item(7, 7).Value = 50.09443433333333
$ws.Cells.Item(7, 8).Value = 150.283303
$ws.Cells.Item(7, 9).Value = 0.4062381250674705
$ws.Cells.Item(7, 10).Value = 0.4062381250674706
$ws.Cells.Item(7, 11).Value = 3
$ws.Cells.Item(7, 13).Value = 7.676201333333334
$ws.Cells.Item(7, 14).Value = 23.028604
$ws.Cells.Item(7, 15).Value = 0.5555816676272323
$ws.Cells.Item(7, 16).Value = 0.5555816676272323
$ws.Cells.Item(7, 17).Value = 384.5349636221125
$ws.Cells.Item(7, 18).Value = 3460.814672599012
$ws.Cells.Item(7, 19).Value = 0.2256984549787454
$ws.Cells.Item(7, 20).Value = 0.2256984549787455

$ws.Cells.Item(8, 5).Value = 3
$ws.Cells.Item(8, 7).Value = 11.30520233333333
$ws.Cells.Item(8, 8).Value = 33.915607
$ws.Cells.Item(8, 9).Value = 0.09167893121303822
$ws.Cells.Item(8, 10).Value = 0.09167893121303823
$ws.Cells.Item(8, 11).Value = 3
$ws.Cells.Item(8, 13).Value = 4.358929333333333
$ws.Cells.Item(8, 14).Value = 13.076788
$ws.Cells.Item(8, 15).Value = 0.315486934607403
$ws.Cells.Item(8, 16).Value = 0.315486934607403
$ws.Cells.Item(8, 17).Value = 49.27857807003511
$ws.Cells.Item(8, 18).Value = 443.507202630316
$ws.Cells.Item(8, 19).Value = 0.02892350497648438
$ws.Cells.Item(8, 20).Value = 0.02892350497648439

$ws.Cells.Item(9, 5).Value = 3
$ws.Cells.Item(9, 7).Value = 11.30520233333333
$ws.Cells.Item(9, 8).Value = 33.915607
$ws.Cells.Item(9, 9).Value = 0.09167893121303822
$ws.Cells.Item(9, 10).Value = 0.09167893121303823
$ws.Cells.Item(9, 11).Value = 3
$ws.Cells.Item(9, 13).Value = 1.781382333333333
$ws.Cells.Item(9, 14).Value = 5.344147
$ws.Cells.Item(9, 15).Value = 0.1289313977653647
$ws.Cells.Item(9, 16).Value = 0.1289313977653647
$ws.Cells.Item(9, 17).Value = 20.13888771135878
$ws.Cells.Item(9, 18).Value = 181.249989402229
$ws.Cells.Item(9, 19).Value = 0.01182029274693174
$ws.Cells.Item(9, 20).Value = 0.01182029274693175

$ws.Cells.Item(10, 5).Value = 3
$ws.Cells.Item(10, 7).Value = 11.30520233333333
$ws.Cells.Item(10, 8).Value = 33.915607
$ws.Cells.Item(10, 9).Value = 0.09167893121303822
$ws.Cells.Item(10, 10).Value = 0.09167893121303823
$ws.Cells.Item(10, 11).Value = 3
$ws.Cells.Item(10, 13).Value = 7.676201333333334
$ws.Cells.Item(10, 14).Value = 23.028604
$ws.Cells.Item(10, 15).Value = 0.5555816676272323
$ws.Cells.Item(10, 16).Value = 0.5555816676272323
$ws.Cells.Item(10, 17).Value = 86.78100922473645
$ws.Cells.Item(10, 18).Value = 781.0290830226281
$ws.Cells.Item(10, 19).Value = 0.05093513348962209
$ws.Cells.Item(10, 20).Value = 0.05093513348962209
